$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new "season record" columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new header cells the same formatting as the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team's season record (90 wins, 72 losses, 0 ties) for every player row
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}
